$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    # Force text storage so numeric-looking strings (e.g. "1.00", "538.66")
    # keep their exact formatting instead of being parsed as numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    # Restore the original cell style so we do not introduce formatting changes
    # that are not part of the intended edit.
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '63.144.29'
Set-TextValue 'E2' '  -1.94%  '
Set-TextValue 'D3' '2.628.22'
Set-TextValue 'E3' '  -1.78%  '
Set-TextValue 'E4' '  -0.06%  '
Set-TextValue 'D5' '605.89'
Set-TextValue 'E5' '  +1.67%  '
Set-TextValue 'D6' '146.56'
Set-TextValue 'E6' '  -0.95%  '
Set-TextValue 'D7' '1.00'
Set-TextValue 'E7' '  +0.00%  '
Set-TextValue 'D8' '0.586'
Set-TextValue 'E8' '  -1.23%  '
Set-TextValue 'D9' '2.627.73'
Set-TextValue 'E9' '  -1.79%  '
Set-TextValue 'E10' '  -0.22%  '
Set-TextValue 'E12' '  +0.03%  '
Set-TextValue 'E13' '  +1.15%  '
Set-TextValue 'D14' '27.25'
Set-TextValue 'E14' '  -2.50%  '
Set-TextValue 'D15' '3.094.52'
Set-TextValue 'E15' '  -1.93%  '
Set-TextValue 'D16' '62.981.44'
Set-TextValue 'E16' '  -2.07%  '
Set-TextValue 'E17' '  -2.07%  '
Set-TextValue 'D18' '2.640.08'
Set-TextValue 'E18' '  -2.69%  '
Set-TextValue 'D19' '11.31'
Set-TextValue 'E19' '  -1.08%  '
Set-TextValue 'E20' '  +1.89%  '
Set-TextValue 'D21' '340.23'
Set-TextValue 'E21' '  -1.83%  '
Set-TextValue 'D22' '6.87'
Set-TextValue 'E22' '  -0.29%  '
Set-TextValue 'D23' '0.999'
Set-TextValue 'E23' '  -0.20%  '
Set-TextValue 'E24' '  -4.56%  '
Set-TextValue 'D25' '66.57'
Set-TextValue 'E25' '  -3.22%  '
Set-TextValue 'D26' '1.63'
Set-TextValue 'E26' '  -2.89%  '
Set-TextValue 'E27' '  -4.58%  '
Set-TextValue 'D28' '8.68'
Set-TextValue 'E28' '  +1.58%  '
Set-TextValue 'E29' '  -2.64%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D30' '1.00'
Set-TextValue 'E30' '  -0.04%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D31' '538.66'
Set-TextValue 'E31' '  +1.65%  '
Set-TextValue 'E32' '  -1.29%  '
Set-TextValue 'D33' '2.04'
Set-TextValue 'E33' '  +1.49%  '
Set-TextValue 'E34' '  -2.36%  '
Set-TextValue 'D35' '0.0₃0804'
Set-TextValue 'E35' '  -2.92%  '
Set-TextValue 'D36' '5.21'
Set-TextValue 'E36' '  +11.15%  '
Set-TextValue 'D37' '169.30'
Set-TextValue 'E37' '  -3.59%  '
Set-TextValue 'E38' '  -0.07%  '
Set-TextValue 'E39' '  -0.11%  '
Set-TextValue 'D40' '19.03'
Set-TextValue 'E40' '  -1.59%  '
Set-TextValue 'D41' '1.88'
Set-TextValue 'E41' '  +6.08%  '
Set-TextValue 'E42' '  +0.01%  '
Set-TextValue 'D43' '169.43'
Set-TextValue 'E43' '  -2.06%  '
Set-TextValue 'D44' '3.75'
Set-TextValue 'E44' '  -1.00%  '
Set-TextValue 'D45' '22.32'
Set-TextValue 'E45' '  +2.37%  '
Set-TextValue 'D46' '0.0568'
Set-TextValue 'E46' '  +2.99%  '
Set-TextValue 'D47' '0.624'
Set-TextValue 'E47' '  -1.86%  '
Set-TextValue 'E48' '  -0.84%  '
Set-TextValue 'D49' '0.0961'
Set-TextValue 'E49' '  -0.31%  '
Set-TextValue 'D50' '18.49'
Set-TextValue 'E50' '  -2.05%  '
Set-TextValue 'E51' '  -0.12%  '
